# Force delete source item data too before updating
#
# - Hide the "data_validation" helper sheet.
# - Update the General sheet's explanatory text (row 8) to match the new
#   "force delete" behaviour, and grow the row to fit the longer text.
# - Update the selections that were left active on the General and
#   Item_policies sheets when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# --- workbook.xml: hide the data_validation sheet -------------------------
$dv = $wb.Worksheets.Item("data_validation")
$dv.Visible = $false

# --- sharedStrings.xml: rewrite the two explanatory cells on General ------
$general = $wb.Worksheets.Item("General")
$general.Range("A8").Value = "In the source IZ delete delete the fields of items when containing defined values"
$general.Range("D8").Value = 'Fields like ''provenance'', ''temp_location'', ''temp_library'', ''in_temp_location'', ''pattern_type'', ''statistics_note_1'', ''statistics_note_2'', ''statistics_note_3'', ''po_line'' can cause error and prevent sometime barcode update (adding prefix "OLD_").'

# --- sheet1.xml: row 8 grows to fit the rewritten text ---------------------
$general.Rows.Item(8).RowHeight = 64.2

# --- sheet4.xml (Item_policies): leftover selection from last save --------
$itemPolicies = $wb.Worksheets.Item("Item_policies")
$itemPolicies.Activate() | Out-Null
$itemPolicies.Range("B2").Select() | Out-Null

# --- sheet1.xml (General): re-activate + leave the new selection ----------
# General must end up as the active/selected tab (tabSelected="1"), so it
# needs to be (re)activated last.
$general.Activate() | Out-Null
$general.Range("A11").Select() | Out-Null
